$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the style of row 50 (for the date-formatted column A) down to rows 51:57
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A51:A57").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 51
$ws.Range("A51").Value = 45794
$ws.Range("B51").Value = "Flowering"
$ws.Range("C51").Value = "Large"
$ws.Range("D51").Value = 55
$ws.Range("E51").Value = 75
$ws.Range("F51").Formula = "=ABS(D51-E51)"
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0.5
$ws.Range("I51").Value = "Yes"
$ws.Range("J51").Value = 2
$ws.Range("K51").Value = "Neutral"
$ws.Range("L51").Value = 6
$ws.Range("M51").Value = 0.52
$ws.Range("N51").Value = 55
$ws.Range("O51").Value = 29.55
$ws.Range("P51").Value = 35
$ws.Range("Q51").Value = 0.76
$ws.Range("R51").Value = 9.9
$ws.Range("S51").Value = 43
$ws.Range("T51").Value = 40

# Row 52
$ws.Range("A52").Value = 45794
$ws.Range("B52").Value = "Nonflowering"
$ws.Range("C52").Value = "Medium"
$ws.Range("D52").Value = 55
$ws.Range("E52").Value = 75
$ws.Range("F52").Formula = "=ABS(D52-E52)"
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0.25
$ws.Range("I52").Value = "Yes"
$ws.Range("J52").Value = 3
$ws.Range("K52").Value = "Bright"
$ws.Range("L52").Value = 6
$ws.Range("M52").Value = 0.52
$ws.Range("N52").Value = 55
$ws.Range("O52").Value = 29.55
$ws.Range("P52").Value = 35
$ws.Range("Q52").Value = 0.76
$ws.Range("R52").Value = 9.9
$ws.Range("S52").Value = 43
$ws.Range("T52").Value = 40

# Row 53
$ws.Range("A53").Value = 45794
$ws.Range("B53").Value = "Nonflowering"
$ws.Range("C53").Value = "Small"
$ws.Range("D53").Value = 55
$ws.Range("E53").Value = 75
$ws.Range("F53").Formula = "=ABS(D53-E53)"
$ws.Range("G53").Value = 0
$ws.Range("H53").Formula = "=1/3"
$ws.Range("I53").Value = "Yes"
$ws.Range("J53").Value = 3
$ws.Range("K53").Value = "Bright"
$ws.Range("L53").Value = 6
$ws.Range("M53").Value = 0.52
$ws.Range("N53").Value = 55
$ws.Range("O53").Value = 29.55
$ws.Range("P53").Value = 35
$ws.Range("Q53").Value = 0.76
$ws.Range("R53").Value = 9.9
$ws.Range("S53").Value = 43
$ws.Range("T53").Value = 40

# Row 54
$ws.Range("A54").Value = 45794
$ws.Range("B54").Value = "Nonflowering"
$ws.Range("C54").Value = "Medium"
$ws.Range("D54").Value = 55
$ws.Range("E54").Value = 75
$ws.Range("F54").Formula = "=ABS(D54-E54)"
$ws.Range("G54").Value = 0
$ws.Range("H54").Formula = "=2/3"
$ws.Range("I54").Value = "Yes"
$ws.Range("J54").Value = 3
$ws.Range("K54").Value = "Neutral"
$ws.Range("L54").Value = 6
$ws.Range("M54").Value = 0.52
$ws.Range("N54").Value = 55
$ws.Range("O54").Value = 29.55
$ws.Range("P54").Value = 35
$ws.Range("Q54").Value = 0.76
$ws.Range("R54").Value = 9.9
$ws.Range("S54").Value = 43
$ws.Range("T54").Value = 40

# Row 55
$ws.Range("A55").Value = 45794
$ws.Range("B55").Value = "Nonflowering"
$ws.Range("C55").Value = "Medium"
$ws.Range("D55").Value = 55
$ws.Range("E55").Value = 75
$ws.Range("F55").Formula = "=ABS(D55-E55)"
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0.1
$ws.Range("I55").Value = "Yes"
$ws.Range("J55").Value = 3
$ws.Range("K55").Value = "Dark"
$ws.Range("L55").Value = 6
$ws.Range("M55").Value = 0.52
$ws.Range("N55").Value = 55
$ws.Range("O55").Value = 29.55
$ws.Range("P55").Value = 35
$ws.Range("Q55").Value = 0.76
$ws.Range("R55").Value = 9.9
$ws.Range("S55").Value = 43
$ws.Range("T55").Value = 40

# Row 56
$ws.Range("A56").Value = 45794
$ws.Range("B56").Value = "Nonflowering"
$ws.Range("C56").Value = "Large"
$ws.Range("D56").Value = 55
$ws.Range("E56").Value = 75
$ws.Range("F56").Formula = "=ABS(D56-E56)"
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 0.1
$ws.Range("I56").Value = "Yes"
$ws.Range("J56").Value = 4
$ws.Range("K56").Value = "Dark"
$ws.Range("L56").Value = 6
$ws.Range("M56").Value = 0.52
$ws.Range("N56").Value = 55
$ws.Range("O56").Value = 29.55
$ws.Range("P56").Value = 35
$ws.Range("Q56").Value = 0.76
$ws.Range("R56").Value = 9.9
$ws.Range("S56").Value = 43
$ws.Range("T56").Value = 40

# Row 57
$ws.Range("A57").Value = 45794
$ws.Range("B57").Value = "Tree"
$ws.Range("C57").Value = "Medium"
$ws.Range("D57").Value = 55
$ws.Range("E57").Value = 75
$ws.Range("F57").Formula = "=ABS(D57-E57)"
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 1
$ws.Range("I57").Value = "Yes"
$ws.Range("J57").Value = 1
$ws.Range("K57").Value = "Dark"
$ws.Range("L57").Value = 6
$ws.Range("M57").Value = 0.52
$ws.Range("N57").Value = 55
$ws.Range("O57").Value = 29.55
$ws.Range("P57").Value = 35
$ws.Range("Q57").Value = 0.76
$ws.Range("R57").Value = 9.9
$ws.Range("S57").Value = 43
$ws.Range("T57").Value = 40
